# Update TPM-derived NATMI ligand-receptor metrics (Ybx1-Notch1) with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 156.920329
$ws.Range("H2").Value = 470.760987
$ws.Range("I2").Value = 0.5508291342957625
$ws.Range("J2").Value = 0.5508291342957624
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 6049.698706363956
$ws.Range("R2").Value = 54447.28835727561
$ws.Range("S2").Value = 0.3171757730309882
$ws.Range("T2").Value = 0.3171757730309882
$ws.Range("G3").Value = 156.920329
$ws.Range("H3").Value = 470.760987
$ws.Range("I3").Value = 0.5508291342957625
$ws.Range("J3").Value = 0.5508291342957624
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 842.7682448724041
$ws.Range("R3").Value = 7584.914203851637
$ws.Range("S3").Value = 0.04418495573542907
$ws.Range("T3").Value = 0.04418495573542907
$ws.Range("G4").Value = 156.920329
$ws.Range("H4").Value = 470.760987
$ws.Range("I4").Value = 0.5508291342957625
$ws.Range("J4").Value = 0.5508291342957624
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 3613.853469557822
$ws.Range("R4").Value = 32524.6812260204
$ws.Range("S4").Value = 0.1894684055293452
$ws.Range("T4").Value = 0.1894684055293452
$ws.Range("I5").Value = 0.2834387340807631
$ws.Range("J5").Value = 0.2834387340807631
$ws.Range("M5").Value = 38.55267666666666
$ws.Range("N5").Value = 115.65803
$ws.Range("O5").Value = 0.5758151725879548
$ws.Range("P5").Value = 0.5758151725879548
$ws.Range("Q5").Value = 3112.977938420243
$ws.Range("R5").Value = 28016.80144578219
$ws.Range("S5").Value = 0.163208323582826
$ws.Range("T5").Value = 0.163208323582826
$ws.Range("I6").Value = 0.2834387340807631
$ws.Range("J6").Value = 0.2834387340807631
$ws.Range("O6").Value = 0.08021535714867321
$ws.Range("P6").Value = 0.08021535714867323
$ws.Range("Q6").Value = 433.661093027516
$ws.Range("R6").Value = 3902.949837247645
$ws.Range("S6").Value = 0.02273613928405623
$ws.Range("T6").Value = 0.02273613928405623
$ws.Range("I7").Value = 0.2834387340807631
$ws.Range("J7").Value = 0.2834387340807631
$ws.Range("M7").Value = 23.02986166666667
$ws.Range("N7").Value = 69.089585
$ws.Range("O7").Value = 0.3439694702633719
$ws.Range("P7").Value = 0.3439694702633719
$ws.Range("Q7").Value = 1859.571305854078
$ws.Range("R7").Value = 16736.1417526867
$ws.Range("S7").Value = 0.09749427121388081
$ws.Range("T7").Value = 0.09749427121388082
$ws.Range("G8").Value = 47.21380733333333
$ws.Range("H8").Value = 141.641422
$ws.Range("I8").Value = 0.1657321316234745
$ws.Range("J8").Value = 0.1657321316234745
$ws.Range("M8").Value = 38.55267666666666
$ws.Range("N8").Value = 115.65803
$ws.Range("O8").Value = 0.5758151725879548
$ws.Range("P8").Value = 0.5758151725879548
$ws.Range("Q8").Value = 1820.218648324295
$ws.Range("R8").Value = 16381.96783491866
$ws.Range("S8").Value = 0.09543107597414059
$ws.Range("T8").Value = 0.09543107597414059
$ws.Range("G9").Value = 47.21380733333333
$ws.Range("H9").Value = 141.641422
$ws.Range("I9").Value = 0.1657321316234745
$ws.Range("J9").Value = 0.1657321316234745
$ws.Range("O9").Value = 0.08021535714867321
$ws.Range("P9").Value = 0.08021535714867323
$ws.Range("Q9").Value = 253.5700619137573
$ws.Range("R9").Value = 2282.130557223816
$ws.Range("S9").Value = 0.01329426212918792
$ws.Range("T9").Value = 0.01329426212918793
$ws.Range("G10").Value = 47.21380733333333
$ws.Range("H10").Value = 141.641422
$ws.Range("I10").Value = 0.1657321316234745
$ws.Range("J10").Value = 0.1657321316234745
$ws.Range("M10").Value = 23.02986166666667
$ws.Range("N10").Value = 69.089585
$ws.Range("O10").Value = 0.3439694702633719
$ws.Range("P10").Value = 0.3439694702633719
$ws.Range("Q10").Value = 1087.327451643319
$ws.Range("R10").Value = 9785.947064789869
$ws.Range("S10").Value = 0.05700679352014593
$ws.Range("T10").Value = 0.05700679352014595
